# Correct the misspelled category label "Scuidal " -> "Suicidal".
# The typo appears twice in column A of Sheet1 (used as a lookup/category key
# for populating the DB), so find every occurrence and fix it rather than
# hard-coding row numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    if ($current -ne $null -and $current.ToString().Trim() -eq "Scuidal") {
        $cell.Value = "Suicidal"
    }
}
